$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.297.80"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "2.477.84"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "521.62"
$ws.Range("E5").Value = "  -2.84%  "
$ws.Range("D6").Value = "132.35"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.559"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("D11").Value = "5.37"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "0.345"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "2.917.81"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "58.240.67"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "22.18"
$ws.Range("E15").Value = "  -3.68%  "
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").Value = "2.477.90"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").Value = "10.87"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "4.20"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "321.35"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  -2.90%  "
$ws.Range("D23").Value = "64.33"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").Value = "0.410"
$ws.Range("E24").Value = "  -2.90%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  -3.14%  "
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "1.21"
$ws.Range("E29").Value = "  +3.06%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "6.38"
$ws.Range("E30").Value = "  -5.15%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.71"
$ws.Range("E31").Value = "  -3.94%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "166.11"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "18.20"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("E36").Value = "  -10.00%  "
$ws.Range("D37").Value = "4.02"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").Value = "0.798"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").Value = "278.68"
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("E41").Value = "  -3.99%  "
$ws.Range("D42").Value = "5.13"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").Value = "0.598"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Value = "126.23"
$ws.Range("E44").Value = "  -4.70%  "
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("E46").Value = "  -2.90%  "
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").Value = "17.29"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "1.745.31"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("E51").Value = "  -1.38%  "
